# The deck ships with two DrawingML themes embedded in the package:
#   ppt/theme/theme1.xml -> "Integral" / "Red Violet" (the theme actually
#                            used by the (only) slide master, i.e. the
#                            presentation's active design)
#   ppt/theme/theme2.xml -> "Office Theme" / "Office" (an unused, orphaned
#                            theme part not referenced by any master/slide)
#
# The authored change swaps the two themes around: the slide master's
# design becomes the (former) "Office Theme" colour palette. Apply this
# through the supported PowerPoint object model by re-pointing every
# slot of the active design's ThemeColorScheme (the 12-colour DrawingML
# clrScheme: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) at the
# "Office Theme" RGB values.
#
# Colour values are passed as OLE/VBA colour longs (0x00BBGGRR), i.e.
# R + G*256 + B*65536 for a given RRGGBB hex colour.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

$scheme.Item(1).RGB  = 0         # dk1      000000
$scheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$scheme.Item(3).RGB  = 6968388   # dk2      44546A
$scheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$scheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$scheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$scheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$scheme.Item(8).RGB  = 49407     # accent4  FFC000
$scheme.Item(9).RGB  = 12874308  # accent5  4472C4
$scheme.Item(10).RGB = 4697456   # accent6  70AD47
$scheme.Item(11).RGB = 12673797  # hlink    0563C1
$scheme.Item(12).RGB = 7491477   # folHlink 954F72
